$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Replace-ParagraphXml($paragraph, $innerXml) {
    # InsertXML only splices correctly in place when the target Range spans an
    # entire paragraph (start through, and including, the paragraph mark), so
    # we always operate on $paragraph.Range as a whole.
    $rng = $paragraph.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# --- Paragraph containing "name = {m:v.name}," ---------------------------
# The run holding "{m" is split into two runs: "{" and "m", so the parser's
# TokenIteratorFieldRewriterSplit sees the opening brace as its own token.
$p3 = $d.Paragraphs(3)
if ($p3.Range.Text -notlike "*{m:v.name}*") {
    throw "Unexpected content for paragraph 3: $($p3.Range.Text)"
}
$p3xml = '<w:p w:rsidP="00F5495F" w:rsidR="00052FB8" w:rsidRDefault="00730F02">' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>name</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>' + $nbsp + '</w:t></w:r>' +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r><w:t>=</w:t></w:r>' +
         '<w:r w:rsidR="00E27251"><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t>{</w:t></w:r>' +
         '<w:r><w:t>m</w:t></w:r>' +
         '<w:r w:rsidR="006F5523"><w:t xml:space="preserve">:v.name}</w:t></w:r>' +
         '<w:r w:rsidR="00052FB8"><w:t>,</w:t></w:r>' +
         '</w:p>'
Replace-ParagraphXml $p3 $p3xml

# --- Paragraph containing "{m:endfor}" ------------------------------------
# The run holding "{m:" is split into two runs: "{" and "m:".
$p4 = $d.Paragraphs(4)
if ($p4.Range.Text -notlike "*{m:endfor}*") {
    throw "Unexpected content for paragraph 4: $($p4.Range.Text)"
}
$p4xml = '<w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="006F5523">' +
         '<w:r><w:t>{</w:t></w:r>' +
         '<w:r><w:t>m:</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">endfor}</w:t></w:r>' +
         '</w:p>'
Replace-ParagraphXml $p4 $p4xml

Write-Host $d.Content.Text
